$wb = $excel.ActiveWorkbook

# --- "2016" sheet updates ---
$ws2016 = $wb.Worksheets.Item("2016")

# Row 3 (EMR): shares and December dividend payment updated
$ws2016.Range("D3").Value = 14.372999999999999
$ws2016.Range("S3").Value = 6.84

# Row 5 (RDS.A / Shell): shares and December dividend payment updated
$ws2016.Range("D5").Value = 40.033000000000001
$ws2016.Range("S5").Value = 13.9

# Row 8 (RSD.A / Archer Daniels): shares and December dividend payment updated
$ws2016.Range("D8").Value = 18.715
$ws2016.Range("S8").Value = 17.260000000000002

# Restore the active selection shown in the sheet to I14
$ws2016.Activate()
$ws2016.Range("I14").Select()
